# Weekly price-sheet update: a new week of "Plátano" price observations is
# inserted at the top of the data block (right after row 233), pushing every
# existing observation down by two rows. The sheet grows from 350 to 352 rows
# (header + 351 data rows) and the two brand-new rows get fresh figures while
# keeping the same "Calidad" (quality) labels that used to live in that slot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right after the header/unchanged block (row 233),
# which pushes the former row 234 down to 236, former 235 down to 237, etc.
$ws.Rows.Item(234).Insert()
$ws.Rows.Item(234).Insert()

# The two new rows reuse the same "shape" (all non-price/date columns) as
# the rows that used to sit there before being pushed down to 236 / 237.
for ($c = 1; $c -le 20; $c++) {
    $ws.Cells.Item(234, $c).Value = $ws.Cells.Item(236, $c).Value()
    $ws.Cells.Item(235, $c).Value = $ws.Cells.Item(237, $c).Value()
}

# New observation for row 234 (Pintón)
$ws.Cells.Item(234, 4).Value  = 44455   # Fecha
$ws.Cells.Item(234, 13).Value = 1050    # Volumen
$ws.Cells.Item(234, 14).Value = 17000   # Precio mínimo
$ws.Cells.Item(234, 15).Value = 17000   # Precio máximo
$ws.Cells.Item(234, 16).Value = 17000   # Precio promedio ponderado
$ws.Cells.Item(234, 19).Value = 850     # Precio $/Kg

# New observation for row 235 (Primera Pintón)
$ws.Cells.Item(235, 4).Value  = 44455   # Fecha
$ws.Cells.Item(235, 13).Value = 500     # Volumen
$ws.Cells.Item(235, 14).Value = 18000   # Precio mínimo
$ws.Cells.Item(235, 15).Value = 18000   # Precio máximo
$ws.Cells.Item(235, 16).Value = 18000   # Precio promedio ponderado
$ws.Cells.Item(235, 19).Value = 900     # Precio $/Kg
